$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Mark the "CU - Consultar historial de pagos/rentas de cliente" task (row 7) as done
$ws.Range("F7").Value = "Hecho"

# Register 4 hours consumed on day 1 (column K) for that task
$ws.Range("K7").Value = 4

# Leave selection on D13 as in the edited workbook
$ws.Range("D13").Select()
